# Cotações atualizadas - 2025-09-30
# Append a new daily quote row (row 26) to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

# Date column keeps the same date/time number format as the rows above it.
$ws.Range("A$row").Value = 45930
$ws.Range("A$row").NumberFormat = $ws.Range("A25").NumberFormat

# Quote columns are stored as text, matching the existing rows' formatting.
$ws.Range("B$row").Value = "21,2348"
$ws.Range("C$row").Value = "15,0941"
$ws.Range("D$row").Value = "15,0177"
$ws.Range("E$row").Value = "15,0177"
